# Adds the new monthly ranking rows (2024-04-30 snapshot) to the "Planilha1"
# tracking sheet, mirroring the layout already used for the prior months.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "quality" strategy rows (62-71) - same 10 tickers/ranks used for the
#    previous quality snapshots, dated 2024-04-30 (serial 45412).
# ---------------------------------------------------------------------------
$quality = @(
    @("CMIN3", 1),
    @("CPFE3", 2),
    @("CSNA3", 3),
    @("ENGI11", 4),
    @("EQTL3", 5),
    @("JBSS3", 6),
    @("MRFG3", 7),
    @("PETR3", 8),
    @("PRIO3", 9),
    @("VBBR3", 10)
)

$r = 62
foreach ($item in $quality) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = "quality"
    $ws.Cells.Item($r, 4).Value = 45412
    $ws.Range("D61").Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $r++
}

# ---------------------------------------------------------------------------
# 2) New "long" strategy rows (72-81).
# ---------------------------------------------------------------------------
$long = @(
    @("BRFS3", 1),
    @("EMBR3", 2),
    @("STBP3", 3),
    @("NTCO3", 4),
    @("PETR4", 5),
    @("SBSP3", 6),
    @("DIRR3", 7),
    @("TEND3", 8),
    @("GMAT3", 9),
    @("MRFG3", 10)
)

# ---------------------------------------------------------------------------
# 3) New "short" strategy rows (82-91).
# ---------------------------------------------------------------------------
$short = @(
    @("MRVE3", 10),
    @("VAMO3", 9),
    @("LWSA3", 8),
    @("ARZZ3", 7),
    @("BEEF3", 6),
    @("MOVI3", 5),
    @("BHIA3", 4),
    @("AZUL4", 3),
    @("PCAR3", 2),
    @("CVCB3", 1)
)

# Build the black Calibri font used to highlight column A for this new batch
# of rows (rows 72-91) once, on the first cell; every following cell reuses
# the same style via a format-only paste so we do not keep allocating new
# style entries for every row.
$ws.Cells.Item(72, 1).Value = $long[0][0]
$ws.Cells.Item(72, 1).Font.Name = "Calibri"
$ws.Cells.Item(72, 1).Font.Color = 0
$ws.Rows.Item(72).RowHeight = 15

$r = 72
foreach ($item in $long) {
    if ($r -ne 72) {
        $ws.Cells.Item($r, 1).Value = $item[0]
        $ws.Range("A72").Copy()
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
        $ws.Rows.Item($r).RowHeight = 15
    }
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = "long"
    $ws.Cells.Item($r, 4).Value = 45412
    $ws.Range("D61").Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $r++
}

foreach ($item in $short) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Range("A72").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 15
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = "short"
    $ws.Cells.Item($r, 4).Value = 45412
    $ws.Range("D61").Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $r++
}

# ---------------------------------------------------------------------------
# 4) Clean up the stray highlight style that used to sit on C59 - it is not
#    part of the data and the refreshed sheet no longer carries it.
# ---------------------------------------------------------------------------
$ws.Range("C59").Style = "Normal"

# ---------------------------------------------------------------------------
# 5) Move the view to the newly added rows, matching where the sheet was
#    left scrolled to after the update.
# ---------------------------------------------------------------------------
$ws.Range("D82:D91").Select()
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1
